$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-26 Sunday", "2025-01-27 Monday"),
    @("194÷6=32, 2", "219÷9=24, 3"),
    @("138÷9=15, 3", "888÷4=222, 0"),
    @("492÷4=123, 0", "459÷6=76, 3"),
    @("292÷4=73, 0", "809÷3=269, 2"),
    @("406÷9=45, 1", "210÷2=105, 0"),
    @("457÷4=114, 1", "850÷8=106, 2"),
    @("449÷6=74, 5", "716÷5=143, 1"),
    @("354÷2=177, 0", "776÷4=194, 0"),
    @("868÷8=108, 4", "969÷6=161, 3"),
    @("747÷8=93, 3", "237÷8=29, 5"),
    @("533÷5=106, 3", "893÷9=99, 2"),
    @("266÷5=53, 1", "197÷5=39, 2"),
    @("880÷8=110, 0", "132÷4=33, 0"),
    @("315÷2=157, 1", "822÷4=205, 2"),
    @("434÷3=144, 2", "240÷2=120, 0"),
    @("164÷9=18, 2", "702÷9=78, 0"),
    @("612÷7=87, 3", "510÷8=63, 6"),
    @("272÷3=90, 2", "419÷8=52, 3"),
    @("167÷5=33, 2", "781÷9=86, 7"),
    @("161÷6=26, 5", "984÷5=196, 4"),
    @("366÷6=61, 0", "535÷6=89, 1"),
    @("362÷7=51, 5", "188÷5=37, 3"),
    @("840÷8=105, 0", "959÷3=319, 2"),
    @("566÷4=141, 2", "449÷5=89, 4"),
    @("602÷3=200, 2", "882÷9=98, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    [void]$range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
